$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,("2018-10", 101.4, 100.6, 100.3)
    ,("2018-11", 101.4, 100.6, 100.6)
    ,("2018-12", 101.2, 100.6, 103.5)
    ,("2018-01", 100.2198, 100.0177, 100.5062)
    ,("2018-02", 100.3, 100, 99.90000000000001)
    ,("2018-03", 100.4, 100, 100)
    ,("2018-04", 101.3, 100, 100.1)
    ,("2018-05", 101.3, 100.3, 99.8)
    ,("2018-06", 101.4, 100.5, 99.8)
    ,("2018-07", 100.6, 100.5, 99.8)
    ,("2018-08", 100.6, 100.6, 99.8)
    ,("2018-09", 100.9, 100.6, 99.8)
    ,("2019-10", 101.4, 102.9, 103)
    ,("2019-11", 101.3, 102.9, 102.6)
    ,("2019-12", 101.3, 102.9, 99.7)
    ,("2019-01", 103, 100.7, 102.1)
    ,("2019-02", 103, 100.7, 102.7)
    ,("2019-03", 103, 100.7, 102.8)
    ,("2019-04", 101.6, 103, 102.9)
    ,("2019-05", 101.7, 103, 103.2)
    ,("2019-06", 101.6, 102.8, 103.2)
    ,("2019-07", 101.9, 102.8, 103.2)
    ,("2019-08", 101.9, 102.8, 103.2)
    ,("2019-09", 101.9, 102.9, 103.2)
    ,("2020-10", 102.7, 100.8, 101.3)
    ,("2020-11", 102.8, 100.8, 101.3)
    ,("2020-12", 103, 100.7, 101.3)
    ,("2020-01", 99.40000000000001, 103, 101.3)
    ,("2020-02", 99.3, 102.9, 101.6)
    ,("2020-03", 103, 103.3, 101.8)
    ,("2020-04", 103, 101, 101.4)
    ,("2020-05", 102.9, 100.7, 101.4)
    ,("2020-06", 102.9, 100.7, 101.3)
    ,("2020-07", 103, 100.8, 101.3)
    ,("2020-08", 103, 100.8, 101.3)
    ,("2020-09", 102.7, 100.7, 101.3)
    ,("2021-10", 100.7, 100.5, 100)
    ,("2021-11", 100.7, 100.6, 100)
    ,("2021-12", 100.7, 100.9, 100)
    ,("2021-01", 99.7, 100.6, 100.7)
    ,("2021-02", 99.7, 100.6, 100)
    ,("2021-03", 100.6, 100.4, 99.7)
    ,("2021-04", 100.6, 100.4, 100)
    ,("2021-05", 100.6, 100.5, 100)
    ,("2021-06", 100.7, 100.7, 100)
    ,("2021-07", 100.7, 100.6, 100)
    ,("2021-08", 100.7, 100.6, 100)
    ,("2021-09", 100.7, 100.6, 100)
    ,("2022-10", 99.8, 100.5, 100.4)
    ,("2022-11", 99.8, 100.4, 100.4)
    ,("2022-12", 99.8, 100.2, 100.4)
    ,("2022-01", 101.4, 101, 100.4)
    ,("2022-02", 101.4, 101, 100.4)
    ,("2022-03", 99.90000000000001, 100.9, 100.4)
    ,("2022-04", 99.90000000000001, 100.9, 100.4)
    ,("2022-05", 99.8, 100.8, 100.4)
    ,("2022-06", 99.8, 100.5, 100.4)
    ,("2022-07", 99.8, 100.5, 100.4)
    ,("2022-08", 99.8, 100.5, 100.4)
    ,("2022-09", 99.8, 100.5, 100.4)
    ,("2023-01", 99.90000000000001, 100.5, 100.4)
    ,("2023-02", 99.90000000000001, 100.5, 100.4)
    ,("2023-03", 100.1, 100.5, 100.8)
    ,("2023-04", 100.1, 100.5, 100.8)
    ,("2023-05", 100.2, 100.5, 100.8)
    ,("2023-06", 100.2, 100.5, 100.8)
    ,("2023-07", 104.4, 101.3, 100.8)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}